$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.554.11"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.629.21"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").Value = "'213.16"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'0.498"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "'0.251"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'19.07"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "1.855.36"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.607.11"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").Value = "'4.12"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "'0.521"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "'63.70"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "26.576.83"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'215.24"
$ws.Range("E19").Value = "  +5.41%  "
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").Value = "'4.31"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'6.16"
$ws.Range("D23").Value = "'9.36"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("E24").Value = "  +4.32%  "
$ws.Range("D25").Value = "'148.08"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'6.87"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D32").Value = "'3.31"
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.220.29"
$ws.Range("E35").Value = "  +5.20%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.38"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("E37").Value = "  +5.52%  "
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "'0.797"
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "'0.499"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").Value = "'2.27"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("D42").Value = "'0.796"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("D43").Value = "'5.35"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").Value = "1.765.53"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "'92.59"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("D47").Value = "'55.07"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "'7.67"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").Value = "'0.408"
$ws.Range("E51").Value = "  -0.26%  "
